$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticket rows for 2024-05-27 (rows 276-286). Force column A to text
# so the date-like "YYYY-MM-DD" strings are not auto-converted to Excel
# date serial numbers (matches the source data, stored as literal text).
$ws.Range("A276:A286").NumberFormat = "@"

# Row 276
$ws.Range("A276").Value = "2024-05-27"
$ws.Range("B276").Value = "10:51:31"
$ws.Range("C276").Value = "-"
$ws.Range("D276").Value = "Cámara no detecta busbar"
$ws.Range("E276").Value = "-"
$ws.Range("F276").Value = "-"
$ws.Range("G276").Value = "-"
$ws.Range("H276").Value = "10:51:33"
$ws.Range("I276").Value = "0:00:02"

# Row 277
$ws.Range("A277").Value = "2024-05-27"
$ws.Range("B277").Value = "10:51:51"
$ws.Range("C277").Value = "-"
$ws.Range("D277").Value = "No detecta presencia power CP"
$ws.Range("E277").Value = "-"
$ws.Range("F277").Value = "-"
$ws.Range("G277").Value = "-"
$ws.Range("H277").Value = "10:51:53"
$ws.Range("I277").Value = "0:00:02"

# Row 278
$ws.Range("A278").Value = "2024-05-27"
$ws.Range("B278").Value = "10:52:01"
$ws.Range("C278").Value = "-"
$ws.Range("D278").Value = "Detección de sealling mal puesto"
$ws.Range("E278").Value = "-"
$ws.Range("F278").Value = "-"
$ws.Range("G278").Value = "-"
$ws.Range("H278").Value = "10:52:03"
$ws.Range("I278").Value = "0:00:02"

# Row 279
$ws.Range("A279").Value = "2024-05-27"
$ws.Range("B279").Value = "10:52:06"
$ws.Range("C279").Value = "-"
$ws.Range("D279").Value = "Detección de sealling mal puesto"
$ws.Range("E279").Value = "-"
$ws.Range("F279").Value = "-"
$ws.Range("G279").Value = "-"
$ws.Range("H279").Value = "10:52:08"
$ws.Range("I279").Value = "0:00:02"

# Row 280
$ws.Range("A280").Value = "2024-05-27"
$ws.Range("B280").Value = "10:54:19"
$ws.Range("C280").Value = "-"
$ws.Range("D280").Value = "Cámara no detecta foams"
$ws.Range("E280").Value = "-"
$ws.Range("F280").Value = "-"
$ws.Range("G280").Value = "-"
$ws.Range("H280").Value = "10:54:21"
$ws.Range("I280").Value = "0:00:02"

# Row 281
$ws.Range("A281").Value = "2024-05-27"
$ws.Range("B281").Value = "10:54:24"
$ws.Range("C281").Value = "-"
$ws.Range("D281").Value = "Cámara no detecta foams"
$ws.Range("E281").Value = "-"
$ws.Range("F281").Value = "-"
$ws.Range("G281").Value = "-"
$ws.Range("H281").Value = "10:54:26"
$ws.Range("I281").Value = "0:00:02"

# Row 282
$ws.Range("A282").Value = "2024-05-27"
$ws.Range("B282").Value = "10:56:29"
$ws.Range("C282").Value = "-"
$ws.Range("D282").Value = "No coloca bien el sealling"
$ws.Range("E282").Value = "-"
$ws.Range("F282").Value = "-"
$ws.Range("G282").Value = "-"
$ws.Range("H282").Value = "10:56:32"
$ws.Range("I282").Value = "0:00:03"

# Row 283
$ws.Range("A283").Value = "2024-05-27"
$ws.Range("B283").Value = "11:12:33"
$ws.Range("C283").Value = "-"
$ws.Range("D283").Value = "No coloca bien el sealling"
$ws.Range("E283").Value = "-"
$ws.Range("F283").Value = "-"
$ws.Range("G283").Value = "-"
$ws.Range("H283").Value = "11:12:36"
$ws.Range("I283").Value = "0:00:03"

# Row 284
$ws.Range("A284").Value = "2024-05-27"
$ws.Range("B284").Value = "11:12:38"
$ws.Range("C284").Value = "-"
$ws.Range("D284").Value = "Cámara no detecta busbar"
$ws.Range("E284").Value = "-"
$ws.Range("F284").Value = "-"
$ws.Range("G284").Value = "-"
$ws.Range("H284").Value = "11:12:40"
$ws.Range("I284").Value = "0:00:02"

# Row 285
$ws.Range("A285").Value = "2024-05-27"
$ws.Range("B285").Value = "11:24:06"
$ws.Range("C285").Value = "-"
$ws.Range("D285").Value = "Cámara no detecta Pcb"
$ws.Range("E285").Value = "-"
$ws.Range("F285").Value = "-"
$ws.Range("G285").Value = "-"
$ws.Range("H285").Value = "11:24:09"
$ws.Range("I285").Value = "0:00:03"

# Row 286
$ws.Range("A286").Value = "2024-05-27"
$ws.Range("B286").Value = "11:24:12"
$ws.Range("C286").Value = "-"
$ws.Range("D286").Value = "Cámara no detecta Power CP"
$ws.Range("E286").Value = "-"
$ws.Range("F286").Value = "-"
$ws.Range("G286").Value = "-"
$ws.Range("H286").Value = "11:24:14"
$ws.Range("I286").Value = "0:00:02"
